# Consolidate the fragmented "sehr gut (1) = sgt," style runs (left over
# from older per-character formatting) into single runs per grade entry.
# Word's Find/Replace keeps the formatting of the first run in the matched
# range and removes the runs (and any inline markup, e.g. proofErr) that
# used to carry the remaining characters.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "sehr gut (1) = sgt,", $true, $false, $false, $false, $false,
    $true, 1, $false, "sehr gut (1) = sgt,", 2) | Out-Null

$d.Content.Find.Execute(
    "gut (2) = gut,", $true, $false, $false, $false, $false,
    $true, 1, $false, "gut (2) = gut,", 2) | Out-Null

$d.Content.Find.Execute(
    "befriedigend (3) = bfr,", $true, $false, $false, $false, $false,
    $true, 1, $false, "befriedigend (3) = bfr,", 2) | Out-Null

$d.Content.Find.Execute(
    "ausreichend (4) = ausr,", $true, $false, $false, $false, $false,
    $true, 1, $false, "ausreichend (4) = ausr,", 2) | Out-Null

$d.Content.Find.Execute(
    " mangelhaft (5) =  mgh,", $true, $false, $false, $false, $false,
    $true, 1, $false, " mangelhaft (5) = mgh,", 2) | Out-Null

$d.Content.Find.Execute(
    " ungenügend (6) = ung", $true, $false, $false, $false, $false,
    $true, 1, $false, " ungenügend (6) = ung", 2) | Out-Null
